# Update the "dSF" (column F) values for the rows whose underlying
# source data was re-pulled, per commit message:
# "repull data, push all data, mean calculation"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    4  = -9
    5  = -3
    9  = -4
    13 = 3
    14 = 1
    15 = 5
    16 = 4
    19 = 4
    20 = -4
    21 = -1
    22 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
